$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.89109567700462367
$ws.Range("T1").Value = 0.99733591644918662
$ws.Range("BO1").Value = 0.80230237146098449
$ws.Range("A2").Value = 0.88531050062956573
$ws.Range("C2").Value = 0.67530291010070598
$ws.Range("D2").Value = 0.96075970327740723
$ws.Range("D3").Value = 0.99458008699367817
$ws.Range("N3").Value = 0.56923778774630429
$ws.Range("F4").Value = 0.95304037942569453
$ws.Range("BG4").Value = 0.68910473394214478
$ws.Range("D5").Value = 0.90645144399747857
$ws.Range("F5").Value = 0.88812135821974258
$ws.Range("G6").Value = 0.77757825191047147
$ws.Range("E7").Value = 0.65318026502646143
$ws.Range("W7").Value = 0.57954532768864686
$ws.Range("G8").Value = 0.78505219462538922
$ws.Range("J8").Value = 0.98667463951463241
$ws.Range("H9").Value = 0.77085412752634319
$ws.Range("J9").Value = 0.89799542238828112
$ws.Range("AL9").Value = 0.95203389046332965
$ws.Range("L10").Value = 0.6434675435793018
$ws.Range("Z10").Value = 0.84210377983060969
$ws.Range("I11").Value = 0.59440384699857507
$ws.Range("BA11").Value = 0.71276085544514378
$ws.Range("N12").Value = 0.78103138439626008
$ws.Range("AB12").Value = 0.87316050407382551
$ws.Range("L13").Value = 0.8183789176202142
$ws.Range("N13").Value = 0.94152123273478794
$ws.Range("O13").Value = 0.98180200849469546
$ws.Range("P14").Value = 0.952223360151076
$ws.Range("AI14").Value = 0.97960015364635633
$ws.Range("Q15").Value = 0.68806104569387183
$ws.Range("R15").Value = 0.78537975449466302
$ws.Range("Q16").Value = 0.89256510802341005
$ws.Range("AC16").Value = 0.7188194214105621
$ws.Range("BM17").Value = 0.95857658839526838
$ws.Range("I18").Value = 0.87454913667839473
$ws.Range("Q18").Value = 0.80345875171023706
$ws.Range("S18").Value = 0.75260658497422717
$ws.Range("S20").Value = 0.6127663893817733
$ws.Range("U20").Value = 0.90750933236109099
$ws.Range("AB20").Value = 0.99217558660185012
$ws.Range("S21").Value = 0.79175368545892444
$ws.Range("V21").Value = 0.69450048230431305
$ws.Range("BE21").Value = 0.88147567821479011
$ws.Range("X22").Value = 0.78397457811843063
$ws.Range("K23").Value = 0.79572621674017818
$ws.Range("U23").Value = 0.9157184995986396
$ws.Range("V23").Value = 0.88604890167777928
$ws.Range("AC23").Value = 0.78657466555481204
$ws.Range("AG24").Value = 0.97040099149278591
$ws.Range("AV24").Value = 0.65618528224333073
$ws.Range("Z25").Value = 0.80971444677970439
$ws.Range("BB25").Value = 0.92646902347231364
$ws.Range("Z27").Value = 0.89534648739715172
$ws.Range("AC27").Value = 0.8234207225900585
$ws.Range("AA28").Value = 0.98802646416572559
$ws.Range("AU28").Value = 0.93887969771792634
$ws.Range("AF30").Value = 0.97377271948247368
$ws.Range("AC31").Value = 0.92111646015200899
$ws.Range("AD31").Value = 0.98698669174061182
$ws.Range("W32").Value = 0.64124658824381098
$ws.Range("AE32").Value = 0.79455270016985069
$ws.Range("AH32").Value = 0.99985049390785963
$ws.Range("AE33").Value = 0.66055691379295112
$ws.Range("AF33").Value = 0.92463913776961626
$ws.Range("W35").Value = 0.59834815792857721
$ws.Range("AG35").Value = 0.97856398805539757
$ws.Range("AH35").Value = 0.50796778821360555
$ws.Range("AK35").Value = 0.89433128149665486
$ws.Range("AI36").Value = 0.83408943367489008
$ws.Range("AK36").Value = 0.66371138584082368
$ws.Range("AN36").Value = 0.73705733960544784
$ws.Range("AH37").Value = 0.95771731859823572
$ws.Range("AL37").Value = 0.82298477674625492
$ws.Range("AM37").Value = 0.88757799654034386
$ws.Range("BP37").Value = 0.77642320503301931
$ws.Range("AM38").Value = 0.95681072424932068
$ws.Range("Y39").Value = 0.93305807452199818
$ws.Range("BC40").Value = 0.82100055233368896
$ws.Range("AM41").Value = 0.91439347249515091
$ws.Range("AN41").Value = 0.80748503476365363
$ws.Range("AP41").Value = 0.83860886290377579
$ws.Range("AQ41").Value = 0.90248436652853881
$ws.Range("BE42").Value = 0.82319517359845618
$ws.Range("BO42").Value = 0.89770962161726997
$ws.Range("AS43").Value = 0.98769164229439921
$ws.Range("AQ44").Value = 0.72214474813424379
$ws.Range("AS44").Value = 0.93296874588443779
$ws.Range("AT44").Value = 0.81048230353688111
$ws.Range("AU45").Value = 0.91399864334854009
$ws.Range("R46").Value = 0.66821789769189177
$ws.Range("AS46").Value = 0.62740359255912137
$ws.Range("AV46").Value = 0.96538484830724403
$ws.Range("AT47").Value = 0.92821103764955737
$ws.Range("AV47").Value = 0.71222310102719377
$ws.Range("AW48").Value = 0.7955544560043597
$ws.Range("AJ49").Value = 0.86373509009565419
$ws.Range("AY50").Value = 0.97588462754527949
$ws.Range("AZ50").Value = 0.80451608327600055
$ws.Range("AW51").Value = 0.9646447581624894
$ws.Range("BA51").Value = 0.65622945552929779
$ws.Range("H52").Value = 0.99027759938293314
$ws.Range("AY52").Value = 0.92966153799387985
$ws.Range("AZ53").Value = 0.88941863278427924
$ws.Range("BB53").Value = 0.57478938541313829
$ws.Range("BC53").Value = 0.78127611591882751
$ws.Range("BC54").Value = 0.98281937228983463
$ws.Range("BD54").Value = 0.67657171570754715
$ws.Range("X55").Value = 0.57254219664857797
$ws.Range("BD55").Value = 0.91286137316452565
$ws.Range("O56").Value = 0.9671643501706193
$ws.Range("AR56").Value = 0.761667947922819
$ws.Range("BC57").Value = 0.63365717561130253
$ws.Range("BF57").Value = 0.97848039929403186
$ws.Range("AI58").Value = 0.71788089105286779
$ws.Range("BL58").Value = 0.76044114532941576
$ws.Range("AZ59").Value = 0.6277115672905702
$ws.Range("BF59").Value = 0.68594877180311065
$ws.Range("A60").Value = 0.92514520478040452
$ws.Range("E61").Value = 0.91463519754082434
$ws.Range("BK61").Value = 0.61944875760642937
$ws.Range("BH62").Value = 0.59308474326995386
$ws.Range("J63").Value = 0.72875502172492057
$ws.Range("BF63").Value = 0.69649545154254233
$ws.Range("BJ63").Value = 0.54719651219681087
$ws.Range("AP64").Value = 0.59666358940221398
$ws.Range("BJ64").Value = 0.98645594442322426
$ws.Range("BK64").Value = 0.88151220951304565
$ws.Range("BM64").Value = 0.96306217144888662
$ws.Range("F65").Value = 0.66350196990168486
$ws.Range("AF65").Value = 0.90158014801775721
$ws.Range("BP66").Value = 0.61932598292227026
$ws.Range("BM67").Value = 0.87901645153230712
$ws.Range("BN67").Value = 0.87275854368957551
$ws.Range("AP68").Value = 0.84860910163086456
